$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price/Volume columns hold text-formatted numbers (e.g. "570.42", "0.110").
# Force Text format before writing so Excel does not auto-convert them to
# real numbers (which would drop formatting like trailing zeros), then restore
# the original General format once the values are in place.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '62.628.44'
$ws.Range('E2').Value = '  -0.86%  '
$ws.Range('D3').Value = '2.453.03'
$ws.Range('E3').Value = '  -0.96%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '570.42'
$ws.Range('E5').Value = '  -1.14%  '
$ws.Range('D6').Value = '145.89'
$ws.Range('E6').Value = '  -0.62%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('D8').Value = '0.529'
$ws.Range('E8').Value = '  -1.88%  '
$ws.Range('D9').Value = '0.110'
$ws.Range('E9').Value = '  -0.78%  '
$ws.Range('E10').Value = '  +0.03%  '
$ws.Range('D11').Value = '5.21'
$ws.Range('E11').Value = '  -1.58%  '
$ws.Range('E12').Value = '  -1.91%  '
$ws.Range('D13').Value = '28.56'
$ws.Range('E13').Value = '  -2.02%  '
$ws.Range('E14').Value = '  -3.17%  '
$ws.Range('D15').Value = '2.885.92'
$ws.Range('D16').Value = '62.504.97'
$ws.Range('E16').Value = '  -0.99%  '
$ws.Range('D17').Value = '2.450.87'
$ws.Range('E17').Value = '  -1.01%  '
$ws.Range('E18').Value = '  -6.15%  '
$ws.Range('D19').Value = '10.74'
$ws.Range('E19').Value = '  -3.25%  '
$ws.Range('B20').Value = 'Polkadot'
$ws.Range('C20').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D20').Value = '4.13'
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '320.70'
$ws.Range('E21').Value = '  -3.02%  '
$ws.Range('D22').Value = '2.19'
$ws.Range('E22').Value = '  -1.73%  '
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').Value = '9.90'
$ws.Range('E24').Value = '  +4.52%  '
$ws.Range('D25').Value = '64.96'
$ws.Range('E25').Value = '  -2.27%  '
$ws.Range('D26').Value = '641.80'
$ws.Range('E26').Value = '  -3.95%  '
$ws.Range('D27').Value = '2.560.94'
$ws.Range('E27').Value = '  -1.31%  '
$ws.Range('D28').Value = '0.0₃0955'
$ws.Range('E28').Value = '  -4.09%  '
$ws.Range('D29').Value = '0.997'
$ws.Range('E29').Value = '  -0.16%  '
$ws.Range('D30').Value = '1.42'
$ws.Range('D31').Value = '7.80'
$ws.Range('E31').Value = '  -3.93%  '
$ws.Range('E32').Value = '  -3.33%  '
$ws.Range('D33').Value = '0.131'
$ws.Range('E33').Value = '  -2.35%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('E35').Value = '  -4.00%  '
$ws.Range('E36').Value = '  -3.24%  '
$ws.Range('D37').Value = '150.36'
$ws.Range('E37').Value = '  -2.20%  '
$ws.Range('D38').Value = '0.364'
$ws.Range('E38').Value = '  -2.41%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D39').Value = '5.32'
$ws.Range('E39').Value = '  -3.46%  '
$ws.Range('B40').Value = 'EthereumClassic'
$ws.Range('C40').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D40').Value = '18.48'
$ws.Range('E40').Value = '  -1.63%  '
$ws.Range('E41').Value = '  -0.56%  '
$ws.Range('E42').Value = '  -3.02%  '
$ws.Range('D43').Value = '0.0₆0310'
$ws.Range('E43').Value = '  +2.00%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').Value = '152.74'
$ws.Range('E45').Value = '  +2.52%  '
$ws.Range('D46').Value = '15.39'
$ws.Range('E46').Value = '  +1.48%  '
$ws.Range('D47').Value = '3.53'
$ws.Range('E47').Value = '  -2.70%  '
$ws.Range('E48').Value = '  -0.73%  '
$ws.Range('D49').Value = '20.11'
$ws.Range('E49').Value = '  -4.03%  '
$ws.Range('E50').Value = '  -2.45%  '
$ws.Range('D51').Value = '0.0905'
$ws.Range('E51').Value = '  -1.46%  '

$dataRange.NumberFormat = "General"
